$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.107.30'
$ws.Range('E2').Value = '  +1.69%  '

$ws.Range('D3').Value = '2.507.39'
$ws.Range('E3').Value = '  +0.55%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.29%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.12%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.529'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.19%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('E9').Value = '  +1.20%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.08%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.11'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +9.39%  '

$ws.Range('E12').Value = '  +0.80%  '

$ws.Range('E13').Value = '  +0.72%  '

$ws.Range('E14').Value = '  +0.72%  '

$ws.Range('D15').Value = '2.901.37'
$ws.Range('E15').Value = '  +0.64%  '

$ws.Range('D16').Value = '2.549.68'
$ws.Range('E16').Value = '  +2.08%  '

$ws.Range('E17').Value = '  -0.14%  '

$ws.Range('D18').Value = '47.955.43'
$ws.Range('E18').Value = '  +1.54%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.15%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.28%  '

$ws.Range('D21').Value = '0.0₃0941'
$ws.Range('E21').Value = '  +0.73%  '

$ws.Range('E22').Value = '  +3.00%  '

$ws.Range('E23').Value = '  +2.47%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '274.82'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +12.11%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.55'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.24%  '

$ws.Range('E26').Value = '  +0.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.49%  '

$ws.Range('E29').Value = '  +0.54%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.140'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.25%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.26%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.48'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.63%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.30'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.51%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.34'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.35%  '

$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.00%  '

$ws.Range('E36').Value = '  -0.18%  '

$ws.Range('E37').Value = '  -0.58%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.63'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.76%  '

$ws.Range('E39').Value = '  +1.30%  '

$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '122.55'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.99%  '

$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.112'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.54%  '

$ws.Range('E42').Value = '  -0.65%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.81'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.50%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0304'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.91%  '

$ws.Range('D45').Value = '2.019.97'
$ws.Range('E45').Value = '  +1.22%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.13'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.72%  '

$ws.Range('E47').Value = '  +4.22%  '

$ws.Range('E48').Value = '  -0.71%  '

$ws.Range('E49').Value = '  -1.79%  '

$ws.Range('E50').Value = '  +1.55%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.53'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.81%  '
